$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new feed-log rows (row 64 and row 65) to Sheet1,
# matching the latest data-lake sync: run_id 63/64, rss_url_id 1/2,
# both captured at 2024-06-16 06:17:11 with a 200 response.

$ws.Range("A64").Value = 63
$ws.Range("B64").Value = 1
$ws.Range("C64").Value = "2024-06-16 06:17:11"
$ws.Range("D64").Value = 200
$ws.Range("E64").Value = 7

$ws.Range("A65").Value = 64
$ws.Range("B65").Value = 2
$ws.Range("C65").Value = "2024-06-16 06:17:11"
$ws.Range("D65").Value = 200
$ws.Range("E65").Value = 0
